$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "设计团队核心职责文档" "设计团队核心责任文档"
Replace-Text "目的" "宗旨"
Replace-Text "：本文档概述了图形设计研究所所有设计团队成员的核心职责。" "：本文档概述了 Graphic Design Institute 所有设计团队成员的核心职责。"
Replace-Text "：与其他设计人员、开发人员和利益干系人协作，创建满足项目要求的高质量设计。" "：与其他设计师、开发人员和利益干系人合作，创造符合项目要求的高质量设计。"
Replace-Text "：创建具有视觉吸引力的设计，这些设计对用户友好、可访问和响应性强。" "：创造具有视觉吸引力、便于用户使用、易于访问且响应速度快的设计。"
Replace-Text "：与团队成员、利益干系人和客户有效沟通，以确保满足项目要求。" "：与团队成员、利益干系人和客户进行有效沟通，确保满足项目要求。"
Replace-Text "：进行研究以确定用户需求、偏好和行为，以告知设计决策。" "：进行研究以确定用户需求、偏好和行为，告知设计决策。"
Replace-Text "：进行可用性测试，以确保设计满足用户需求，可供所有用户访问。" "：进行可用性测试，以确保设计满足用户需求，便于所有用户使用。"
Replace-Text "：随时了解最新的设计趋势、工具和技术，以提高设计质量和效率。" "：随时了解最新设计趋势、工具和技术，以提高设计质量和效率。"
Replace-Text "领导：" "领导力："
Replace-Text " 领导设计团队，为初级设计师提供指导。" "领导设计团队，为初级设计师提供指导。"
